# Realestate Update resale numbers 2025-01-11 22:45
# Append a new data row (row 20) to the CityResaleNum sheet with the
# latest resale numbers snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 20

# Text columns. A (date-like) and D (numeric-like "01") need a leading
# apostrophe so Excel stores them as literal text instead of auto-converting
# to a date serial / number; Style is reset to "Normal" afterwards so the
# quote-prefix formatting doesn't leave a stray number format on the cell.
$ws.Cells.Item($row, 1).Value = "'2025-01-11"
$ws.Cells.Item($row, 1).Style = "Normal"
$ws.Cells.Item($row, 2).Value = "22:45:34"
$ws.Cells.Item($row, 3).Value = "Saturday"
$ws.Cells.Item($row, 4).Value = "'01"
$ws.Cells.Item($row, 4).Style = "Normal"

# Numeric columns (city resale numbers, -1 = no data).
$ws.Cells.Item($row, 5).Value = 127246
$ws.Cells.Item($row, 6).Value = 143585
$ws.Cells.Item($row, 7).Value = 169570
$ws.Cells.Item($row, 8).Value = 159667
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 142945
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 193024
$ws.Cells.Item($row, 14).Value = 115471
$ws.Cells.Item($row, 15).Value = 45886
$ws.Cells.Item($row, 16).Value = 28511
$ws.Cells.Item($row, 17).Value = 65336
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 49292
$ws.Cells.Item($row, 20).Value = -1
